$d = $word.ActiveDocument

# Update the date header (unique text in the document)
$d.Content.Find.Execute("2023-12-16 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-12-17 Sunday", 2)

# The worksheet is a single 20-row x 5-column table; only rows 1, 5, 10,
# 15 and 20 (1-based) carry answer text. Several values repeat across the
# table (e.g. "53x84=4452" appears twice with different replacements), so
# a global Find/Replace would be ambiguous -- address each cell directly
# by table/row/column instead.

$t = $d.Tables.Item(1)

$cellUpdates = @(
    @{Row=1;  Col=1; Text="66×72=4752"},
    @{Row=1;  Col=2; Text="22×90=1980"},
    @{Row=1;  Col=3; Text="89×36=3204"},
    @{Row=1;  Col=4; Text="72×96=6912"},
    @{Row=1;  Col=5; Text="79×69=5451"},

    @{Row=5;  Col=1; Text="17×74=1258"},
    @{Row=5;  Col=2; Text="19×37=703"},
    @{Row=5;  Col=3; Text="13×35=455"},
    @{Row=5;  Col=4; Text="41×23=943"},
    @{Row=5;  Col=5; Text="62×77=4774"},

    @{Row=10; Col=1; Text="45×23=1035"},
    @{Row=10; Col=2; Text="15×44=660"},
    @{Row=10; Col=3; Text="13×93=1209"},
    @{Row=10; Col=4; Text="80×18=1440"},
    @{Row=10; Col=5; Text="60×59=3540"},

    @{Row=15; Col=1; Text="64×34=2176"},
    @{Row=15; Col=2; Text="73×81=5913"},
    @{Row=15; Col=3; Text="24×20=480"},
    @{Row=15; Col=4; Text="69×98=6762"},
    @{Row=15; Col=5; Text="36×64=2304"},

    @{Row=20; Col=1; Text="51×34=1734"},
    @{Row=20; Col=2; Text="24×40=960"},
    @{Row=20; Col=3; Text="64×79=5056"},
    @{Row=20; Col=4; Text="53×56=2968"},
    @{Row=20; Col=5; Text="18×66=1188"}
)

foreach ($u in $cellUpdates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $r = $cell.Range
    # Trim the trailing cell-end / paragraph marks so only the visible text
    # is replaced, keeping the existing run formatting.
    $r.MoveEnd(12, -1) | Out-Null
    $r.Text = $u.Text
}
